# Auto-generated: refresh cryptos list price/volume(1h) cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.769.06"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "'3.525.96"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'607.79"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'196.82"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E9").Value = "  -6.34%  "
$ws.Range("D10").Value = "'0.651"
$ws.Range("D11").Value = "'53.86"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "'4.082.78"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'596.88"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "'12.84"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "'19.14"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "'69.935.95"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'3.524.83"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'18.24"
$ws.Range("E22").Value = "  +6.30%  "
$ws.Range("E23").Value = "  +4.98%  "
$ws.Range("D24").Value = "'4.69"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "'102.29"
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").Value = "'3.19"
$ws.Range("E26").Value = "  +5.30%  "
$ws.Range("D27").Value = "'10.88"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'33.52"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("E31").Value = "  +8.44%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'63.16"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +10.52%  "
$ws.Range("D36").Value = "'3.731.46"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").Value = "'3.08"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "'0.394"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'488.10"
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  +11.57%  "
